$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: new timesheet entry - "Front Design"
$ws.Range("A10").Value = 43978
$ws.Range("B10").Value = 0.7284722222222223
$ws.Range("C10").Value = 0.80763888888888891
$ws.Range("E10").Value = "Front Design"

# Row 11: new timesheet entry - "Front Design Home page register navbar"
$ws.Range("A11").Value = 43979
$ws.Range("B11").Value = 0.42083333333333334
$ws.Range("C11").Value = 0.60486111111111118
$ws.Range("E11").Value = "Front Design Home page register navbar"

# Move the active selection to A12 (matches author's final cursor position)
$ws.Range("A12").Select()
